# Insert a new data row at row 69 (pushing rows 69..183 down to 70..184)
# and populate it with the new record's values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(69).Insert()

$ws.Range("A69").Value = 11
$ws.Range("B69").Value = "Vega Monumental Concepción"
$ws.Range("C69").Value = "Bíobío"
$ws.Range("D69").Value = 44546
$ws.Range("E69").Value = 8
$ws.Range("F69").Value = 100112008
$ws.Range("G69").Value = "Coliflor"
$ws.Range("H69").Value = "Sin especificar"
$ws.Range("I69").Value = "Primera"
$ws.Range("J69").Value = 2500
$ws.Range("K69").Value = 600
$ws.Range("L69").Value = 650
$ws.Range("M69").Value = 630
$ws.Range("N69").Value = "$/unidad"
$ws.Range("O69").Value = "Región Metropolitana"
$ws.Range("P69").Value = 630
$ws.Range("Q69").Value = 1
$ws.Range("R69").Value = "Hortaliza"
